$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column D holds the base "Stardust earned" values; rows 2-201 currently contain 250.
# Update them to 25. Columns E (=D*2) and F (=D*3) contain formulas and will
# recalculate automatically to 50 and 75 respectively.
$ws.Range("D2:D201").Value = 25

# Update the active selection to match the saved view state.
$ws.Range("E8").Select()

$excel.Calculate()
